$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 3065717
$ws.Range("C3").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("P3").Value = "Fjölåsberget, S om, Vrm"
$ws.Range("S3").Value = 25
$ws.Range("X3").Value = "S-Arv-0575"
$ws.Range("AW3").Value = "Värmland Floraväktarna"
$ws.Range("AX3").Value = "Per Larsson"
$ws.Range("AY3").Value = "Floraväkteri Sverige"
